# Update NATMI TPM ligand-receptor pairwise stats (Wnt7b-Fzd4) with
# refreshed values, and expand the sending/target cluster grid to include
# the new "Resolving-Mac" cluster (rows 2-9 cover all FAPs/MuSCs sending x
# ECs/FAPs/MuSCs/Resolving-Mac target combinations).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt7b"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05619466666666667
$ws.Range("H2").Value = 0.168584
$ws.Range("I2").Value = 0.04986276087265156
$ws.Range("J2").Value = 0.07297477932340853
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.524618
$ws.Range("N2").Value = 58.573854
$ws.Range("O2").Value = 0.4154885426712971
$ws.Range("P2").Value = 0.4539723485554654
$ws.Range("Q2").Value = 1.097179400304
$ws.Range("R2").Value = 9.874614602736001
$ws.Range("S2").Value = 0.02071740584854537
$ws.Range("T2").Value = 0.03312853195476459

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt7b"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05619466666666667
$ws.Range("H3").Value = 0.168584
$ws.Range("I3").Value = 0.04986276087265156
$ws.Range("J3").Value = 0.07297477932340853
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.24435933333334
$ws.Range("N3").Value = 45.73307800000001
$ws.Range("O3").Value = 0.324403614112412
$ws.Range("P3").Value = 0.3544508583357054
$ws.Range("Q3").Value = 0.8566516912835558
$ws.Range("R3").Value = 7.709865221552001
$ws.Range("S3").Value = 0.01617565983671113
$ws.Range("T3").Value = 0.02586597316804084

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt7b"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05619466666666667
$ws.Range("H4").Value = 0.168584
$ws.Range("I4").Value = 0.04986276087265156
$ws.Range("J4").Value = 0.07297477932340853
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 11.9507005
$ws.Range("N4").Value = 23.901401
$ws.Range("O4").Value = 0.2543137660693869
$ws.Range("P4").Value = 0.1852460510065796
$ws.Range("Q4").Value = 0.6715656310306667
$ws.Range("R4").Value = 4.029393786184
$ws.Range("S4").Value = 0.01268078650414129
$ws.Range("T4").Value = 0.01351828969273803

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt7b"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05619466666666667
$ws.Range("H5").Value = 0.168584
$ws.Range("I5").Value = 0.04986276087265156
$ws.Range("J5").Value = 0.07297477932340853
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.272275
$ws.Range("N5").Value = 0.816825
$ws.Range("O5").Value = 0.005794077146903843
$ws.Range("P5").Value = 0.006330742102249548
$ws.Range("Q5").Value = 0.01530040286666667
$ws.Range("R5").Value = 0.1377036258
$ws.Range("S5").Value = 0.0002889086832537616
$ws.Range("T5").Value = 0.0004619845078650721

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Wnt7b"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.070792
$ws.Range("H6").Value = 2.141584
$ws.Range("I6").Value = 0.9501372391273485
$ws.Range("J6").Value = 0.9270252206765914
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.524618
$ws.Range("N6").Value = 58.573854
$ws.Range("O6").Value = 0.4154885426712971
$ws.Range("P6").Value = 0.4539723485554654
$ws.Range("Q6").Value = 20.906804757456
$ws.Range("R6").Value = 125.440828544736
$ws.Range("S6").Value = 0.3947711368227518
$ws.Range("T6").Value = 0.4208438166007008

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Wnt7b"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.070792
$ws.Range("H7").Value = 2.141584
$ws.Range("I7").Value = 0.9501372391273485
$ws.Range("J7").Value = 0.9270252206765914
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 15.24435933333334
$ws.Range("N7").Value = 45.73307800000001
$ws.Range("O7").Value = 0.324403614112412
$ws.Range("P7").Value = 0.3544508583357054
$ws.Range("Q7").Value = 16.32353801925867
$ws.Range("R7").Value = 97.94122811555201
$ws.Range("S7").Value = 0.3082279542757009
$ws.Range("T7").Value = 0.3285848851676645

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Wnt7b"
$ws.Range("C8").Value = "Fzd4"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.070792
$ws.Range("H8").Value = 2.141584
$ws.Range("I8").Value = 0.9501372391273485
$ws.Range("J8").Value = 0.9270252206765914
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 11.9507005
$ws.Range("N8").Value = 23.901401
$ws.Range("O8").Value = 0.2543137660693869
$ws.Range("P8").Value = 0.1852460510065796
$ws.Range("Q8").Value = 12.796714489796
$ws.Range("R8").Value = 51.186857959184
$ws.Range("S8").Value = 0.2416329795652457
$ws.Range("T8").Value = 0.1717277613138416

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Wnt7b"
$ws.Range("C9").Value = "Fzd4"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.070792
$ws.Range("H9").Value = 2.141584
$ws.Range("I9").Value = 0.9501372391273485
$ws.Range("J9").Value = 0.9270252206765914
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.272275
$ws.Range("N9").Value = 0.816825
$ws.Range("O9").Value = 0.005794077146903843
$ws.Range("P9").Value = 0.006330742102249548
$ws.Range("Q9").Value = 0.2915498918
$ws.Range("R9").Value = 1.7492993508
$ws.Range("S9").Value = 0.005505168463650082
$ws.Range("T9").Value = 0.005868757594384475
